$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D price cells to remain plain text (matches source
# workbook which stores these as inline strings, not numbers),
# otherwise Excel COM auto-converts numeric-looking strings like
# '0.605' or '1.00' into real numbers and drops formatting.
$ws.Range('D2:D51').NumberFormat = '@'

$ws.Range('D2').Value = '62.858.95'
$ws.Range('E2').Value = '  -1.38%  '
$ws.Range('D3').Value = '3.215.01'
$ws.Range('E3').Value = '  -1.77%  '
$ws.Range('E4').Value = '  -0.30%  '
$ws.Range('D5').Value = '525.57'
$ws.Range('E5').Value = '  +2.29%  '
$ws.Range('D6').Value = '172.17'
$ws.Range('E6').Value = '  -4.04%  '
$ws.Range('E7').Value = '  +1.15%  '
$ws.Range('E8').Value = '  -0.22%  '
$ws.Range('D9').Value = '3.218.55'
$ws.Range('E9').Value = '  -1.21%  '
$ws.Range('D10').Value = '0.605'
$ws.Range('E10').Value = '  -0.60%  '
$ws.Range('D11').Value = '53.14'
$ws.Range('E11').Value = '  -7.42%  '
$ws.Range('E12').Value = '  +2.81%  '
$ws.Range('E13').Value = '  +1.13%  '
$ws.Range('D14').Value = '9.08'
$ws.Range('E14').Value = '  +1.45%  '
$ws.Range('D15').Value = '3.727.58'
$ws.Range('E15').Value = '  -2.62%  '
$ws.Range('E16').Value = '  -4.45%  '
$ws.Range('D17').Value = '3.214.56'
$ws.Range('E17').Value = '  -2.35%  '
$ws.Range('D18').Value = '17.19'
$ws.Range('E18').Value = '  +1.11%  '
$ws.Range('D19').Value = '62.698.31'
$ws.Range('E19').Value = '  -1.49%  '
$ws.Range('D20').Value = '11.03'
$ws.Range('E20').Value = '  +3.23%  '
$ws.Range('D21').Value = '0.965'
$ws.Range('E21').Value = '  +3.33%  '
$ws.Range('D22').Value = '364.80'
$ws.Range('E22').Value = '  -0.64%  '
$ws.Range('D23').Value = '3.75'
$ws.Range('E23').Value = '  +3.51%  '
$ws.Range('E24').Value = '  +2.40%  '
$ws.Range('D25').Value = '11.01'
$ws.Range('E25').Value = '  +3.36%  '
$ws.Range('D26').Value = '3.91'
$ws.Range('E26').Value = '  +5.45%  '
$ws.Range('D27').Value = '6.12'
$ws.Range('E27').Value = '  +3.00%  '
$ws.Range('D28').Value = '2.64'
$ws.Range('E28').Value = '  +1.53%  '
$ws.Range('D29').Value = '11.26'
$ws.Range('E29').Value = '  +2.29%  '
$ws.Range('D30').Value = '8.14'
$ws.Range('E30').Value = '  -0.56%  '
$ws.Range('D31').Value = '28.42'
$ws.Range('E31').Value = '  +1.08%  '
$ws.Range('D32').Value = '628.03'
$ws.Range('E32').Value = '  -2.18%  '
$ws.Range('D33').Value = '6.44'
$ws.Range('E33').Value = '  -1.88%  '
$ws.Range('D34').Value = '11.24'
$ws.Range('E34').Value = '  +2.66%  '
$ws.Range('E35').Value = '  +4.11%  '
$ws.Range('D36').Value = '56.77'
$ws.Range('E36').Value = '  -2.34%  '
$ws.Range('E37').Value = '  +0.08%  '
$ws.Range('D38').Value = '36.66'
$ws.Range('E38').Value = '  +2.93%  '
$ws.Range('D39').Value = '0.374'
$ws.Range('E39').Value = '  +1.90%  '
$ws.Range('D40').Value = '1.00'
$ws.Range('E41').Value = '  +15.39%  '
$ws.Range('E42').Value = '  +1.79%  '
$ws.Range('D43').Value = '2.869.61'
$ws.Range('E43').Value = '  +3.35%  '
$ws.Range('D44').Value = '2.52'
$ws.Range('E44').Value = '  +11.21%  '
$ws.Range('E45').Value = '  +3.99%  '
$ws.Range('E46').Value = '  +11.10%  '
$ws.Range('E47').Value = '  +2.96%  '
$ws.Range('D48').Value = '2.58'
$ws.Range('E48').Value = '  -2.50%  '
$ws.Range('D49').Value = '2.97'
$ws.Range('E49').Value = '  +8.12%  '
$ws.Range('E50').Value = '  +1.61%  '
$ws.Range('D51').Value = '135.18'
$ws.Range('E51').Value = '  +1.35%  '
